$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.868.76'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.72%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.635.85'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.14'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.89%  '
$ws.Range('E6').Value = '  -0.94%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  -1.28%  '
$ws.Range('E9').Value = '  -2.73%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0616'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.07%  '
$ws.Range('E11').Value = '  +1.22%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.868.27'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.81%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.637.41'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.64%  '
$ws.Range('E14').Value = '  -0.62%  '
$ws.Range('E15').Value = '  -0.21%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.42'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.53%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.898.69'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.56%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '230.89'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.78%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0724'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.13%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.55'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.71%  '
$ws.Range('E21').Value = '  -0.06%  '
$ws.Range('E22').Value = '  -0.68%  '
$ws.Range('E23').Value = '  -3.26%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.06'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.90'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.67%  '
$ws.Range('E26').Value = '  +0.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.66'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.85%  '
$ws.Range('E28').Value = '  -1.08%  '
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('E30').Value = '  -0.65%  '
$ws.Range('E31').Value = '  -0.49%  '
$ws.Range('E32').Value = '  +1.77%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.408.64'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.65%  '
$ws.Range('E34').Value = '  -0.07%  '
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('E36').Value = '  +9.15%  '
$ws.Range('E37').Value = '  +1.34%  '
$ws.Range('E38').Value = '  +0.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.564'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.94%  '
$ws.Range('E40').Value = '  -2.46%  '
$ws.Range('E41').Value = '  +0.31%  '
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('E43').Value = '  -3.67%  '
$ws.Range('E44').Value = '  +2.07%  '
$ws.Range('E45').Value = '  +0.14%  '
$ws.Range('E46').Value = '  -1.04%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.777.96'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.88%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '87.92'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0₆0105'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.73%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.1000'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.42%  '
$ws.Range('E51').Value = '  -0.38%  '
